$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "maximum" values (column E) for rows 4-10 from 5 to 4
$ws.Range("E4:E10").Value = 4

# Update the active cell selection to match the final save state
$ws.Range("J21").Select()
